$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.835.14"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.681.49"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.71"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.62"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.95"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.23%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.47%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.48"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000202"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +6.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.163.23"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.650.72"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.682.01"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.63"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.64%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.97"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.26"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.57%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +7.78%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.12%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.20%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.168"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.18"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.15"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "532.84"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.77"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.57"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.428"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.70"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.98"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.36"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "166.36"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0620"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.10"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0264"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.650"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.37"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0989"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.39%  "
